$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 97.5
$ws.Range("I4").Value = 97.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 97.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 16.5
$ws.Range("N4").ClearContents()
# Row 17
$ws.Range("H17").Value = 168391.5
$ws.Range("J17").Value = 168391.5
$ws.Range("L17").Value = 505174.5
$ws.Range("N17").Value = -505510.5
# Row 32
$ws.Range("H32").Value = 17093.125
$ws.Range("J32").Value = 17093.125
$ws.Range("L32").Value = 17093.125
$ws.Range("N32").Value = -17745.125
# Row 40
$ws.Range("H40").Value = 2926460.8
$ws.Range("I40").Value = 2018.2
$ws.Range("J40").Value = 8550389
$ws.Range("K40").Value = 2018.2
$ws.Range("L40").Value = 8550389
$ws.Range("M40").Value = -1843.2
$ws.Range("N40").Value = -8550739
# Row 51
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
# Row 106
$ws.Range("H106").Value = 2770.6667
$ws.Range("I106").Value = 2515.2856
$ws.Range("K106").Value = 2515.2856
$ws.Range("M106").Value = -1884.2856
# Row 132
$ws.Range("H132").Value = 1174.0646
$ws.Range("I132").Value = 1224.6072
$ws.Range("K132").Value = 3673.8216
$ws.Range("M132").Value = -1143.8216
# Row 138
$ws.Range("H138").Value = 2507.3794
$ws.Range("I138").Value = 2308.4736
$ws.Range("J138").Value = 2885.3
$ws.Range("K138").Value = 6925.4208
$ws.Range("L138").Value = 8655.900000000001
$ws.Range("M138").Value = -1785.4208
$ws.Range("N138").Value = -18935.9
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8812.723
$ws.Range("I32").Value = 8243.637000000001
$ws.Range("K32").Value = 8243.637000000001
$ws.Range("M32").Value = -7956.637000000001
# Row 61
$ws.Range("H61").Value = 250012860
$ws.Range("I61").Value = 250012860
$ws.Range("K61").Value = 250012860
$ws.Range("M61").Value = -250012648
# Row 136
$ws.Range("H136").Value = 250012860
$ws.Range("I136").Value = 250012860
$ws.Range("K136").Value = 750038580
$ws.Range("M136").Value = -750036030
# Row 141
$ws.Range("H141").Value = 93996.8
$ws.Range("J141").Value = 93996.8
$ws.Range("L141").Value = 93996.8
$ws.Range("N141").Value = -104356.8
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 112078.555
$ws.Range("I107").Value = 1088.5
$ws.Range("K107").Value = 1088.5
$ws.Range("M107").Value = 831.5
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11840.808
$ws.Range("I31").Value = 6045.5835
$ws.Range("K31").Value = 6045.5835
$ws.Range("M31").Value = -5750.5835
# Row 34
$ws.Range("H34").Value = 11840.808
$ws.Range("I34").Value = 6045.5835
$ws.Range("K34").Value = 6045.5835
$ws.Range("M34").Value = -5843.5835
# Row 62
$ws.Range("H62").Value = 4789.3335
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 4789.3335
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 86
$ws.Range("H86").Value = 4035.4443
$ws.Range("I86").Value = 3865.8
$ws.Range("K86").Value = 3865.8
$ws.Range("M86").Value = -2742.8
# Row 89
$ws.Range("H89").Value = 4035.4443
$ws.Range("I89").Value = 3865.8
$ws.Range("K89").Value = 19329
$ws.Range("M89").Value = -13713
# Row 122
$ws.Range("H122").Value = 2970.7693
$ws.Range("I122").Value = 3464
$ws.Range("J122").Value = 1326.6666
$ws.Range("K122").Value = 10392
$ws.Range("L122").Value = 3979.9998
$ws.Range("M122").Value = -7942
$ws.Range("N122").Value = -8879.9998
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 318.5
$ws.Range("I12").Value = 285.5
$ws.Range("J12").Value = 343.25
$ws.Range("K12").Value = 856.5
$ws.Range("L12").Value = 1029.75
$ws.Range("M12").Value = -683.5
$ws.Range("N12").Value = -1375.75
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
# Row 69
$ws.Range("H69").Value = 709.25
$ws.Range("I69").Value = 709.25
$ws.Range("K69").Value = 2127.75
$ws.Range("M69").Value = -1316.75
# Row 72
$ws.Range("H72").Value = 709.25
$ws.Range("I72").Value = 709.25
$ws.Range("K72").Value = 6383.25
$ws.Range("M72").Value = -2327.25
# Row 80
$ws.Range("H80").Value = 10000
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31872
# Row 83
$ws.Range("H83").Value = 10000
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99360
# Row 107
$ws.Range("H107").Value = 2006.5
$ws.Range("I107").Value = 590.2857
$ws.Range("J107").Value = 3422.7144
$ws.Range("K107").Value = 1770.8571
$ws.Range("L107").Value = 10268.1432
$ws.Range("M107").Value = 149.1428999999998
$ws.Range("N107").Value = -14108.1432
# Row 122
$ws.Range("H122").Value = 522.3077
$ws.Range("I122").Value = 324.1
$ws.Range("J122").Value = 1183
$ws.Range("K122").Value = 2916.9
$ws.Range("L122").Value = 10647
$ws.Range("M122").Value = -466.9000000000001
$ws.Range("N122").Value = -15547
# Row 131
$ws.Range("H131").Value = 2156.7144
$ws.Range("I131").Value = 2750
$ws.Range("J131").Value = 1919.4
$ws.Range("K131").Value = 8250
$ws.Range("L131").Value = 5758.200000000001
$ws.Range("M131").Value = -3210
$ws.Range("N131").Value = -15838.2
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4155.3335
$ws.Range("I46").Value = 1479.6
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 1479.6
$ws.Range("L46").Value = 7500
$ws.Range("M46").Value = -1291.6
$ws.Range("N46").Value = -7876
# Row 61
$ws.Range("H61").Value = 2857.3635
$ws.Range("I61").Value = 2857.3635
$ws.Range("K61").Value = 2857.3635
$ws.Range("M61").Value = -2655.3635
# Row 82
$ws.Range("H82").Value = 878.3
$ws.Range("I82").Value = 898.7778
$ws.Range("J82").Value = 694
$ws.Range("K82").Value = 898.7778
$ws.Range("L82").Value = 694
$ws.Range("M82").Value = -537.7778
$ws.Range("N82").Value = -1416
# Row 85
$ws.Range("H85").Value = 878.3
$ws.Range("I85").Value = 898.7778
$ws.Range("J85").Value = 694
$ws.Range("K85").Value = 898.7778
$ws.Range("L85").Value = 694
$ws.Range("M85").Value = 349.2222
$ws.Range("N85").Value = -3190
# Row 93
$ws.Range("H93").Value = 1375.8
$ws.Range("I93").Value = 1375.8
$ws.Range("K93").Value = 1375.8
$ws.Range("M93").Value = -127.8
# Row 100
$ws.Range("H100").Value = 12477746
$ws.Range("I100").Value = 19960692
$ws.Range("K100").Value = 19960692
$ws.Range("M100").Value = -19960151
# Row 113
$ws.Range("H113").Value = 2857.3635
$ws.Range("I113").Value = 2857.3635
$ws.Range("K113").Value = 2857.3635
$ws.Range("M113").Value = -687.3634999999999
# Row 122
$ws.Range("H122").Value = 8184.1
$ws.Range("I122").Value = 8184.1
$ws.Range("K122").Value = 24552.3
$ws.Range("M122").Value = -22102.3
# Row 140
$ws.Range("H140").Value = 79625.336
$ws.Range("J140").Value = 79625.336
$ws.Range("L140").Value = 79625.336
